$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Shot diameter" column data (header, unit, and all values)
# while keeping the column's cell styling intact.
$ws.Range("D1:D27").ClearContents()
